$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the data range so numeric-looking strings (e.g. "1.00", "0.898")
# are preserved exactly as literal text instead of being normalized as numbers.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

# Apply the updated cryptocurrency values from the latest data pull.
$ws.Range("D2").Value = "68.996.99"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").Value = "3.364.17"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "586.03"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("D6").Value = "179.03"
$ws.Range("E6").Value = "  +2.37%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E8").Value = "  +1.03%  "
$ws.Range("E9").Value = "  +4.45%  "
$ws.Range("D10").Value = "0.585"
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("E11").Value = "  +6.25%  "
$ws.Range("D12").Value = "0.0000275"
$ws.Range("E12").Value = "  +2.47%  "
$ws.Range("D13").Value = "695.80"
$ws.Range("E13").Value = "  +5.48%  "
$ws.Range("D14").Value = "3.923.87"
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("D15").Value = "8.51"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").Value = "69.052.17"
$ws.Range("E16").Value = "  +1.84%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "0.120"
$ws.Range("E17").Value = "  +1.39%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.368.42"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("D19").Value = "17.60"
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("D20").Value = "11.26"
$ws.Range("E20").Value = "  +2.95%  "
$ws.Range("D21").Value = "0.898"
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("E22").Value = "  +2.55%  "
$ws.Range("D23").Value = "17.08"
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("D24").Value = "101.43"
$ws.Range("E24").Value = "  +3.40%  "
$ws.Range("D25").Value = "3.92"
$ws.Range("E25").Value = "  +2.15%  "
$ws.Range("D26").Value = "2.71"
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D27").Value = "9.57"
$ws.Range("E27").Value = "  +3.56%  "
$ws.Range("D28").Value = "33.51"
$ws.Range("E28").Value = "  +0.84%  "
$ws.Range("D29").Value = "8.61"
$ws.Range("E29").Value = "  +2.41%  "
$ws.Range("D30").Value = "7.05"
$ws.Range("E30").Value = "  -2.35%  "
$ws.Range("D31").Value = "11.12"
$ws.Range("E31").Value = "  +1.85%  "
$ws.Range("D32").Value = "554.30"
$ws.Range("E32").Value = "  -2.33%  "
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("E34").Value = "  +9.84%  "
$ws.Range("D35").Value = "57.69"
$ws.Range("E35").Value = "  +2.54%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").Value = "3.711.09"
$ws.Range("E37").Value = "  +0.83%  "
$ws.Range("D38").Value = "0.142"
$ws.Range("E38").Value = "  +9.03%  "
$ws.Range("D39").Value = "34.85"
$ws.Range("D40").Value = "3.21"
$ws.Range("E40").Value = "  +3.56%  "
$ws.Range("D41").Value = "2.64"
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.0$([char]0x2083)0680"
$ws.Range("E42").Value = "  +3.02%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").Value = "0.339"
$ws.Range("E43").Value = "  +2.25%  "
$ws.Range("D44").Value = "0.0416"
$ws.Range("E44").Value = "  +2.87%  "
$ws.Range("E45").Value = "  -1.72%  "
$ws.Range("E46").Value = "  +2.55%  "
$ws.Range("E47").Value = "  +1.35%  "
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("E49").Value = "  -1.69%  "
$ws.Range("D50").Value = "131.77"
$ws.Range("E50").Value = "  +3.51%  "
$ws.Range("E51").Value = "  -1.62%  "

# Restore the default cell style (the NumberFormat="@" trick above would otherwise
# leave a stray "Text" style applied to the range).
$dataRange.Style = "Normal"

